# Update "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the
# per-language handback report sheets, as part of regenerating
# the handback status report.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-11 16:14:34"
$wsZhCn.Range("H3").Value = "2016-03-11 16:14:57"
$wsZhCn.Range("E5").Value = "2016-03-11 16:14:34"
$wsZhCn.Range("H5").Value = "2016-03-11 16:14:57"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-11 16:14:37"
$wsDeDe.Range("H3").Value = "2016-03-11 16:15:05"
$wsDeDe.Range("E5").Value = "2016-03-11 16:14:37"
$wsDeDe.Range("H5").Value = "2016-03-11 16:15:05"
